$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '63.251.32'
$ws.Range('E2').Value = '  -2.00%  '

$ws.Range('D3').Value = '3.163.00'
$ws.Range('E3').Value = '  +0.10%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '590.54'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -2.38%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '138.83'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -4.30%  '

$ws.Range('E7').Value = '  -0.44%  '

$ws.Range('D8').Value = '3.155.45'
$ws.Range('E8').Value = '  +0.14%  '

$ws.Range('E9').Value = '  -1.36%  '

$ws.Range('E10').Value = '  -2.68%  '

$ws.Range('E11').Value = '  -2.16%  '

$ws.Range('E12').Value = '  -2.66%  '

$ws.Range('E13').Value = '  -3.97%  '

$ws.Range('E14').Value = '  -3.88%  '

$ws.Range('D15').Value = '3.678.52'
$ws.Range('E15').Value = '  -0.05%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.120'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.30%  '

$ws.Range('D17').Value = '3.154.95'
$ws.Range('E17').Value = '  -0.13%  '

$ws.Range('D18').Value = '63.206.99'
$ws.Range('E18').Value = '  -2.03%  '

$ws.Range('E19').Value = '  -3.09%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '475.99'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.59%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.12'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -4.13%  '

$ws.Range('E22').Value = '  -1.99%  '

$ws.Range('E23').Value = '  +0.11%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '84.52'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -4.00%  '

$ws.Range('E25').Value = '  -3.73%  '

$ws.Range('E26').Value = '  +0.04%  '

$ws.Range('E27').Value = '  -1.95%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.19'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.10%  '

$ws.Range('E29').Value = '  -5.22%  '

$ws.Range('E30').Value = '  +1.42%  '

$ws.Range('E31').Value = '  +0.06%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '26.94'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.85%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.106'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -5.61%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.54'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -5.66%  '

$ws.Range('E35').Value = '  -3.17%  '

$ws.Range('E36').Value = '  -4.28%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '52.55'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.84%  '

$ws.Range('D38').Value = '0.0₃0701'
$ws.Range('E38').Value = '  -8.23%  '

$ws.Range('E39').Value = '  -2.10%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '422.37'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -4.93%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.76'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -9.65%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '8.29'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.27%  '

$ws.Range('D43').Value = '2.930.94'
$ws.Range('E43').Value = '  +1.67%  '

$ws.Range('E44').Value = '  -6.17%  '

$ws.Range('E45').Value = '  +0.09%  '

$ws.Range('E46').Value = '  -4.81%  '

$ws.Range('E48').Value = '  -2.32%  '

$ws.Range('E49').Value = '  -0.61%  '

$ws.Range('E50').Value = '  -9.59%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '120.94'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.72%  '
